$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.035834837302759
$ws.Range("D2").Value = 1.038271095357433
$ws.Range("E2").Value = 1.049068466986634
$ws.Range("F2").Value = 1.055904597844594
$ws.Range("I2").Value = 1.038812449299044
$ws.Range("J2").Value = 1.040946445527476
$ws.Range("K2").Value = 1.04105960431029
$ws.Range("L2").Value = 1.051826554784042
$ws.Range("M2").Value = 1.058643773686053
$ws.Range("N2").Value = 1.017595811412561
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.036697336242432
$ws.Range("D3").Value = 1.038904091156046
$ws.Range("E3").Value = 1.049951018649166
$ws.Range("F3").Value = 1.056928557161377
$ws.Range("I3").Value = 1.039024252891043
$ws.Range("J3").Value = 1.041453096681774
$ws.Range("K3").Value = 1.041503047006685
$ws.Range("L3").Value = 1.052521096479768
$ws.Range("M3").Value = 1.05948073034828
$ws.Range("N3").Value = 1.017764911107185
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.037255838866624
$ws.Range("D4").Value = 1.039313911060692
$ws.Range("E4").Value = 1.050522911557451
$ws.Range("F4").Value = 1.057592274803294
$ws.Range("I4").Value = 1.039160100893113
$ws.Range("J4").Value = 1.041780688481957
$ws.Range("K4").Value = 1.041789524003802
$ws.Range("L4").Value = 1.05297069570864
$ws.Range("M4").Value = 1.060022833858758
$ws.Range("N4").Value = 1.017874213203426
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.03749072978861
$ws.Range("D5").Value = 1.039486252727134
$ws.Range("E5").Value = 1.050763530877962
$ws.Range("F5").Value = 1.057871574983335
$ws.Range("I5").Value = 1.039216922698905
$ws.Range("J5").Value = 1.041918348295808
$ws.Range("K5").Value = 1.041909847735413
$ws.Range("L5").Value = 1.053159750112453
$ws.Range("M5").Value = 1.060250861323758
$ws.Range("N5").Value = 1.017920135473626
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.037530174624588
$ws.Range("D6").Value = 1.039515192737311
$ws.Range("E6").Value = 1.050803943344372
$ws.Range("F6").Value = 1.057918486696542
$ws.Range("I6").Value = 1.039226446381403
$ws.Range("J6").Value = 1.041941458448725
$ws.Range("K6").Value = 1.041930044068053
$ws.Range("L6").Value = 1.053191495657802
$ws.Range("M6").Value = 1.060289155556553
$ws.Range("N6").Value = 1.017927844352184
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.037258977113368
$ws.Range("D7").Value = 1.039316213690639
$ws.Range("E7").Value = 1.050526125958094
$ws.Range("F7").Value = 1.057596005753937
$ws.Range("I7").Value = 1.039160861284073
$ws.Range("J7").Value = 1.041782528135314
$ws.Range("K7").Value = 1.041791132213647
$ws.Range("L7").Value = 1.052973221695863
$ws.Range("M7").Value = 1.060025880274934
$ws.Range("N7").Value = 1.017874826930858
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.036126237573933
$ws.Range("D8").Value = 1.038484970863701
$ws.Range("E8").Value = 1.049366558162935
$ws.Range("F8").Value = 1.056250411233662
$ws.Range("I8").Value = 1.038884278016197
$ws.Range("J8").Value = 1.041117720841101
$ws.Range("K8").Value = 1.041209562474688
$ws.Range("L8").Value = 1.052061239496111
$ws.Range("M8").Value = 1.058926515216976
$ws.Range("N8").Value = 1.017652983378218
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.034133388626731
$ws.Range("D9").Value = 1.037022037794994
$ws.Range("E9").Value = 1.047329619827133
$ws.Range("F9").Value = 1.053888153376426
$ws.Range("I9").Value = 1.038387718519064
$ws.Range("J9").Value = 1.039944411894535
$ws.Range("K9").Value = 1.040181285154605
$ws.Range("L9").Value = 1.050455676530674
$ws.Range("M9").Value = 1.056993461026155
$ws.Range("N9").Value = 1.017261190003993
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.03280703869134
$ws.Range("D10").Value = 1.036048064624691
$ws.Range("E10").Value = 1.045976030078141
$ws.Range("F10").Value = 1.052319350497983
$ws.Range("I10").Value = 1.038050543585087
$ws.Range("J10").Value = 1.039161034370784
$ws.Range("K10").Value = 1.039493491478634
$ws.Range("L10").Value = 1.049386360073614
$ws.Range("M10").Value = 1.055707637620354
$ws.Range("N10").Value = 1.01699942803566
$ws.Range("B11").Value = 1.019999999999999
$ws.Range("C11").Value = 1.032233255917391
$ws.Range("D11").Value = 1.035626654383046
$ws.Range("E11").Value = 1.045390965335228
$ws.Range("F11").Value = 1.051641490651969
$ws.Range("I11").Value = 1.037903096203696
$ws.Range("J11").Value = 1.038821558141291
$ws.Range("K11").Value = 1.039195142544764
$ws.Range("L11").Value = 1.04892360061746
$ws.Range("M11").Value = 1.055151560942705
$ws.Range("N11").Value = 1.016885952691454
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.032020209138235
$ws.Range("D12").Value = 1.035470174430657
$ws.Range("E12").Value = 1.045173805225265
$ws.Range("F12").Value = 1.051389921377396
$ws.Range("I12").Value = 1.037848110713524
$ws.Range("J12").Value = 1.038695422207127
$ws.Range("K12").Value = 1.03908424390864
$ws.Range("L12").Value = 1.048751751660552
$ws.Range("M12").Value = 1.054945114742727
$ws.Range("N12").Value = 1.016843783643143
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.032065904678982
$ws.Range("D13").Value = 1.035503737598779
$ws.Range("E13").Value = 1.045220379588252
$ws.Range("F13").Value = 1.051443873951063
$ws.Range("I13").Value = 1.037859915095777
$ws.Range("J13").Value = 1.038722480557583
$ws.Range("K13").Value = 1.039108035580953
$ws.Range("L13").Value = 1.048788611992694
$ws.Range("M13").Value = 1.054989393361729
$ws.Range("N13").Value = 1.016852829912277
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.032215643730527
$ws.Range("D14").Value = 1.035613718660928
$ws.Range("E14").Value = 1.045373011561131
$ws.Range("F14").Value = 1.051620691400988
$ws.Range("I14").Value = 1.037898555508294
$ws.Range("J14").Value = 1.03881113250248
$ws.Range("K14").Value = 1.039185977225521
$ws.Range("L14").Value = 1.048909394695729
$ws.Range("M14").Value = 1.055134493877164
$ws.Range("N14").Value = 1.016882467375319
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.032307913764104
$ws.Range("D15").Value = 1.035681488369951
$ws.Range("E15").Value = 1.04546707426772
$ws.Range("F15").Value = 1.051729663403407
$ws.Range("I15").Value = 1.037922334410583
$ws.Range("J15").Value = 1.038865748707802
$ws.Range("K15").Value = 1.039233989266188
$ws.Range("L15").Value = 1.048983818309314
$ws.Range("M15").Value = 1.05522390910318
$ws.Range("N15").Value = 1.016900725454643
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.032845130152307
$ws.Range("D16").Value = 1.036076039244041
$ws.Range("E16").Value = 1.04601488111399
$ws.Range("F16").Value = 1.052364368326282
$ws.Range("I16").Value = 1.038060298720413
$ws.Range("J16").Value = 1.039183558716988
$ws.Range("K16").Value = 1.039513280841693
$ws.Range("L16").Value = 1.04941707750948
$ws.Range("M16").Value = 1.055744557314372
$ws.Range("N16").Value = 1.017006956299675
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.033182256041133
$ws.Range("D17").Value = 1.036323619082878
$ws.Range("E17").Value = 1.046358787627695
$ws.Range("F17").Value = 1.052762889038925
$ws.Range("I17").Value = 1.038146452757599
$ws.Range("J17").Value = 1.03938284133925
$ws.Range("K17").Value = 1.039688331969554
$ws.Range("L17").Value = 1.049688920452812
$ws.Range("M17").Value = 1.056071332705468
$ws.Range("N17").Value = 1.017073557383154
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.033378947453236
$ws.Range("D18").Value = 1.036468059622046
$ws.Range("E18").Value = 1.04655948355116
$ws.Range("F18").Value = 1.052995478481833
$ws.Range("I18").Value = 1.038196565156351
$ws.Range("J18").Value = 1.039499053538497
$ws.Range("K18").Value = 1.039790385150257
$ws.Range("L18").Value = 1.049847507004814
$ws.Range("M18").Value = 1.056262002074614
$ws.Range("N18").Value = 1.017112392033648
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.033446022853876
$ws.Range("D19").Value = 1.036517315394525
$ws.Range("E19").Value = 1.046627932774078
$ws.Range("F19").Value = 1.053074809014985
$ws.Range("I19").Value = 1.03821362845846
$ws.Range("J19").Value = 1.039538674457543
$ws.Range("K19").Value = 1.039825173935188
$ws.Range("L19").Value = 1.049901585152195
$ws.Range("M19").Value = 1.056327026696032
$ws.Range("N19").Value = 1.017125631491234
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.033146080279167
$ws.Range("D20").Value = 1.036297052849655
$ws.Range("E20").Value = 1.046321879248571
$ws.Range("F20").Value = 1.052720117145014
$ws.Range("I20").Value = 1.038137223699638
$ws.Range("J20").Value = 1.039361462880115
$ws.Range("K20").Value = 1.03966955592661
$ws.Range("L20").Value = 1.049659751666325
$ws.Range("M20").Value = 1.056036265878691
$ws.Range("N20").Value = 1.017066413021315
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.032171547046973
$ws.Range("D21").Value = 1.035581330571098
$ws.Range("E21").Value = 1.045328060848917
$ws.Range("F21").Value = 1.051568617038533
$ws.Range("I21").Value = 1.03788718285412
$ws.Range("J21").Value = 1.0387850277924
$ws.Range("K21").Value = 1.03916302751417
$ws.Range("L21").Value = 1.048873826072185
$ws.Range("M21").Value = 1.055091762471097
$ws.Range("N21").Value = 1.016873740409566
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.031559292577608
$ws.Range("D22").Value = 1.035131620923314
$ws.Range("E22").Value = 1.04470412833808
$ws.Range("F22").Value = 1.050845885500647
$ws.Range("I22").Value = 1.037728717260352
$ws.Range("J22").Value = 1.038422373020809
$ws.Range("K22").Value = 1.038844098993289
$ws.Range("L22").Value = 1.048379918429864
$ws.Range("M22").Value = 1.054498526350826
$ws.Range("N22").Value = 1.016752488288041
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.031883814982935
$ws.Range("D23").Value = 1.035369992205632
$ws.Range("E23").Value = 1.045034799023237
$ws.Range("F23").Value = 1.051228899049677
$ws.Range("I23").Value = 1.037812841614484
$ws.Range("J23").Value = 1.038614644372421
$ws.Range("K23").Value = 1.039013211727501
$ws.Range("L23").Value = 1.048641725541334
$ws.Range("M23").Value = 1.054812953734162
$ws.Range("N23").Value = 1.01681677674547
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03316242638956
$ws.Range("D24").Value = 1.03630905689005
$ws.Range("E24").Value = 1.046338556244252
$ws.Range("F24").Value = 1.052739443493259
$ws.Range("I24").Value = 1.038141394345671
$ws.Range("J24").Value = 1.039371122965959
$ws.Range("K24").Value = 1.03967804017046
$ws.Range("L24").Value = 1.049672931707758
$ws.Range("M24").Value = 1.056052110861963
$ws.Range("N24").Value = 1.017069641289985
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.034648202837311
$ws.Range("D25").Value = 1.037400015756226
$ws.Range("E25").Value = 1.047855453529854
$ws.Range("F25").Value = 1.054497796012523
$ws.Range("I25").Value = 1.038517175355496
$ws.Range("J25").Value = 1.040247951360665
$ws.Range("K25").Value = 1.040447525120595
$ws.Range("L25").Value = 1.050870571766758
$ws.Range("M25").Value = 1.057492700032961
$ws.Range("N25").Value = 1.017362579498867
